$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, [string]$val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "63.782.54"
Set-TextValue "E2" "  +3.33%  "
Set-TextValue "D3" "3.078.59"
Set-TextValue "E3" "  +2.82%  "
Set-TextValue "E4" "  +0.16%  "
Set-TextValue "D5" "553.03"
Set-TextValue "E5" "  +3.96%  "
Set-TextValue "D6" "138.50"
Set-TextValue "E6" "  +4.87%  "
Set-TextValue "D7" "1.00"
Set-TextValue "E7" "  +0.03%  "
Set-TextValue "D8" "3.074.78"
Set-TextValue "E8" "  +2.88%  "
Set-TextValue "E9" "  +2.65%  "
Set-TextValue "D10" "0.151"
Set-TextValue "E10" "  +0.73%  "
Set-TextValue "D11" "6.20"
Set-TextValue "E11" "  +0.78%  "
Set-TextValue "E12" "  +2.81%  "
Set-TextValue "D13" "0.0000227"
Set-TextValue "E13" "  +3.86%  "
Set-TextValue "D14" "34.92"
Set-TextValue "E14" "  +4.51%  "
Set-TextValue "D15" "3.579.65"
Set-TextValue "E15" "  +3.40%  "
Set-TextValue "D16" "63.809.40"
Set-TextValue "E16" "  +3.53%  "
Set-TextValue "D17" "3.076.63"
Set-TextValue "E17" "  +3.29%  "
Set-TextValue "E18" "  -1.09%  "
Set-TextValue "D19" "6.73"
Set-TextValue "E19" "  +3.48%  "
Set-TextValue "D20" "485.95"
Set-TextValue "E20" "  +5.86%  "
Set-TextValue "D21" "13.53"
Set-TextValue "E21" "  +2.83%  "
Set-TextValue "D22" "0.684"
Set-TextValue "E22" "  +1.56%  "
Set-TextValue "D23" "7.19"
Set-TextValue "E23" "  +5.00%  "
Set-TextValue "D24" "81.47"
Set-TextValue "E24" "  +4.93%  "
Set-TextValue "D25" "12.55"
Set-TextValue "E25" "  +6.07%  "
Set-TextValue "D26" "1.00"
Set-TextValue "E26" "  +0.34%  "
Set-TextValue "E27" "  +3.95%  "
Set-TextValue "D28" "8.01"
Set-TextValue "E28" "  +4.54%  "
Set-TextValue "D29" "2.00"
Set-TextValue "E29" "  +9.08%  "
Set-TextValue "D30" "0.998"
Set-TextValue "E30" "  +0.14%  "
Set-TextValue "D31" "26.02"
Set-TextValue "E31" "  +1.98%  "
Set-TextValue "E32" "  +1.75%  "
Set-TextValue "B33" "Stacks"
Set-TextValue "C33" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D33" "2.43"
Set-TextValue "E33" "  +8.55%  "
Set-TextValue "B34" "NEARProtocol"
Set-TextValue "C34" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D34" "5.80"
Set-TextValue "E34" "  +7.79%  "
Set-TextValue "D35" "55.68"
Set-TextValue "E35" "  +0.24%  "
Set-TextValue "D36" "5.98"
Set-TextValue "E36" "  +3.31%  "
Set-TextValue "D37" "471.92"
Set-TextValue "E37" "  +4.02%  "
Set-TextValue "D38" "3.187.66"
Set-TextValue "E38" "  +1.65%  "
Set-TextValue "E39" "  +5.11%  "
Set-TextValue "D40" "0.0397"
Set-TextValue "E40" "  +3.36%  "
Set-TextValue "D41" "0.120"
Set-TextValue "E41" "  +2.50%  "
Set-TextValue "D42" "8.22"
Set-TextValue "E42" "  +2.96%  "
Set-TextValue "D43" "2.57"
Set-TextValue "E43" "  +6.49%  "
Set-TextValue "D44" "27.98"
Set-TextValue "E44" "  +10.41%  "
Set-TextValue "E45" "  +3.27%  "
Set-TextValue "E46" "  -0.10%  "
Set-TextValue "E47" "  +5.47%  "
Set-TextValue "E48" "  +2.32%  "
Set-TextValue "D49" "0.0₃0515"
Set-TextValue "E49" "  +1.94%  "
Set-TextValue "D50" "116.63"
Set-TextValue "E50" "  -3.40%  "
Set-TextValue "E51" "  +5.52%  "
